$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4422.936084921746
$ws.Range("C3").Value = 4333.814608369932
$ws.Range("C4").Value = 4152.536524199848
$ws.Range("C5").Value = 4152.536524199848
$ws.Range("C6").Value = 4152.536524199848
$ws.Range("C7").Value = 4152.536524199848
$ws.Range("C8").Value = 4152.536524199848
$ws.Range("C9").Value = 4117.643147798251
$ws.Range("C10").Value = 4117.643147798251
$ws.Range("C11").Value = 4059.245376128709
$ws.Range("C12").Value = 4059.245376128709
